# Insert a new data row before row 80 (shifts existing rows 80:115 down to 81:116)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(80).Insert()

# Columns that stay constant across all data rows in this block.
$ws.Cells.Item(80, 1).Value = 1
$ws.Cells.Item(80, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(80, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(80, 5).Value = 15
$ws.Cells.Item(80, 6).Value = "Fruta"
$ws.Cells.Item(80, 7).Value = 100102
$ws.Cells.Item(80, 8).Value = "Cítricos"
$ws.Cells.Item(80, 9).Value = 100102005
$ws.Cells.Item(80, 10).Value = "Naranja"
$ws.Cells.Item(80, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(80, 20).Value = 1

# New record-specific values.
$ws.Cells.Item(80, 4).Value = 44917
$ws.Cells.Item(80, 4).NumberFormat = $ws.Cells.Item(81, 4).NumberFormat
$ws.Cells.Item(80, 11).Value = "Valencia"
$ws.Cells.Item(80, 12).Value = "Segunda"
$ws.Cells.Item(80, 13).Value = 300
$ws.Cells.Item(80, 14).Value = 900
$ws.Cells.Item(80, 15).Value = 1000
$ws.Cells.Item(80, 16).Value = 967
$ws.Cells.Item(80, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(80, 19).Value = 967
